# Apply cryptos list price/volume refresh + three row reorders (rows 31-33, 38-39, 47-48)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.698.44"
$ws.Range("E2").Value = "  -3.95%  "

$ws.Range("D3").Value = "3.095.74"
$ws.Range("E3").Value = "  -5.04%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'607.24"
$ws.Range("E5").Value = "  -1.22%  "

$ws.Range("D6").Value = "'144.52"
$ws.Range("E6").Value = "  -8.20%  "

$ws.Range("D8").Value = "3.090.66"
$ws.Range("E8").Value = "  -5.17%  "

$ws.Range("E9").Value = "  -5.22%  "

$ws.Range("E10").Value = "  -8.14%  "

$ws.Range("D11").Value = "'5.17"
$ws.Range("E11").Value = "  -10.79%  "

$ws.Range("E12").Value = "  -6.03%  "

$ws.Range("E13").Value = "  -8.98%  "

$ws.Range("E14").Value = "  -10.48%  "

$ws.Range("D15").Value = "3.597.81"
$ws.Range("E15").Value = "  -5.10%  "

$ws.Range("E16").Value = "  +0.82%  "

$ws.Range("D17").Value = "63.703.28"
$ws.Range("E17").Value = "  -4.03%  "

$ws.Range("D18").Value = "3.088.93"
$ws.Range("E18").Value = "  -5.23%  "

$ws.Range("E19").Value = "  -9.06%  "

$ws.Range("D20").Value = "'473.05"
$ws.Range("E20").Value = "  -6.37%  "

$ws.Range("D21").Value = "'14.51"
$ws.Range("E21").Value = "  -6.12%  "

$ws.Range("E23").Value = "  -5.86%  "

$ws.Range("D24").Value = "'13.46"
$ws.Range("E24").Value = "  -8.00%  "

$ws.Range("D25").Value = "'82.80"
$ws.Range("E25").Value = "  -5.01%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -9.33%  "

$ws.Range("D28").Value = "'8.28"
$ws.Range("E28").Value = "  -9.99%  "

$ws.Range("E29").Value = "  -11.09%  "

$ws.Range("D30").Value = "'6.65"
$ws.Range("E30").Value = "  -5.52%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.111"
$ws.Range("E31").Value = "  -13.89%  "

$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.71"
$ws.Range("E33").Value = "  -5.60%  "

$ws.Range("E34").Value = "  -7.33%  "

$ws.Range("E35").Value = "  -4.61%  "

$ws.Range("D36").Value = "'5.88"
$ws.Range("E36").Value = "  -9.07%  "

$ws.Range("D37").Value = "'51.96"

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0725"
$ws.Range("E38").Value = "  -7.89%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'454.60"
$ws.Range("E39").Value = "  -8.20%  "

$ws.Range("D40").Value = "'2.88"
$ws.Range("E40").Value = "  -14.88%  "

$ws.Range("D41").Value = "'0.0390"
$ws.Range("E41").Value = "  -7.59%  "

$ws.Range("D43").Value = "'8.28"
$ws.Range("E43").Value = "  -6.29%  "

$ws.Range("D44").Value = "2.812.18"
$ws.Range("E44").Value = "  -6.26%  "

$ws.Range("E45").Value = "  -10.55%  "

$ws.Range("D46").Value = "'2.22"
$ws.Range("E46").Value = "  -12.19%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  -5.50%  "

$ws.Range("D49").Value = "'25.82"
$ws.Range("E49").Value = "  -10.73%  "

$ws.Range("E50").Value = "  -5.80%  "

$ws.Range("D51").Value = "'117.70"
$ws.Range("E51").Value = "  -2.09%  "
